$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 50 (2025-02-13, abs_activity) with recalculated values ---
$ws.Range("C50").Value = 9.964676051248118
$ws.Range("D50").Value = 9.229307393718706
$ws.Range("F50").Value = 19.19398344496683

# --- Append new rows 54-69 (2025-02-14 .. 2025-02-17) ---
$newRows = @(
    @("2025-02-14", "abs_activity", 10, 10, 0, 20),
    @("2025-02-14", "rel_activity", 10, 10, 0, 20),
    @("2025-02-14", "abs_sleep", 9.533333333333333, 8.533333333333333, 0, 18.06666666666667),
    @("2025-02-14", "rel_sleep", 0, 0, 0, 0),
    @("2025-02-15", "abs_activity", 10, 10, 0, 20),
    @("2025-02-15", "rel_activity", 10, 10, 0, 20),
    @("2025-02-15", "abs_sleep", 10, 9.800000000000001, 0, 19.8),
    @("2025-02-15", "rel_sleep", 10, 7.986606075101652, 0, 17.98660607510165),
    @("2025-02-16", "abs_activity", 10, 10, 0, 20),
    @("2025-02-16", "rel_activity", 0, 10, 0, 10),
    @("2025-02-16", "abs_sleep", 10, 10, 0, 20),
    @("2025-02-16", "rel_sleep", 10, 8.704137766084667, 0, 18.70413776608467),
    @("2025-02-17", "abs_activity", 10, 9.328620175690039, 0, 19.32862017569004),
    @("2025-02-17", "rel_activity", 0, 0, 0, 0),
    @("2025-02-17", "abs_sleep", 0, 8.566666666666666, 0, 8.566666666666666),
    @("2025-02-17", "rel_sleep", 0, 0, 0, 0)
)

$startRow = 54
$endRow = $startRow + $newRows.Count - 1

# Format column A as text first so date-like strings ("2025-02-14") are stored
# as plain text instead of being auto-converted to date serial numbers.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove the temporary text formatting so the new cells end up with the
# default (unstyled) appearance, matching the rest of the data rows.
$ws.Range("A$startRow`:A$endRow").ClearFormats()
